# realistic costs and binary variable for battery
#
# The "bat" (battery) device table had a bogus all-zero placeholder entry
# in the first data row. That row is removed (the real cost/capacity rows
# shift up) and the "Number" column is renumbered 1..7 to match.
# Also refreshes the active-sheet/selection bookkeeping so "bat" becomes
# the active tab (it was "hp_geo" before).

$wb = $excel.ActiveWorkbook

$bat = $wb.Worksheets.Item("bat")
$hpgeo = $wb.Worksheets.Item("hp_geo")

# Drop the placeholder battery row (all zero cost/capacity); everything
# below shifts up one row.
$bat.Rows("2:2").Delete()

# Renumber the "Number" column back to a clean 1..7 sequence.
$bat.Range("A2").Value = 1
$bat.Range("A3").Value = 2
$bat.Range("A4").Value = 3
$bat.Range("A5").Value = 4
$bat.Range("A6").Value = 5
$bat.Range("A7").Value = 6
$bat.Range("A8").Value = 7

# Update stored selections to match what was left active on each sheet.
$hpgeo.Range("E3").Select()
$bat.Range("E15").Select()

# "bat" becomes the active / selected tab.
$bat.Activate()
